$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.044.32'
$ws.Range("E2").Value = '  -1.59%  '
$ws.Range("D3").Value = '3.684.90'
$ws.Range("E3").Value = '  -2.37%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = '595.58'
$ws.Range("E5").Value = '  +0.57%  '
$ws.Range("D6").Value = '165.75'
$ws.Range("E6").Value = '  -2.90%  '
$ws.Range("D7").Value = '3.685.20'
$ws.Range("E7").Value = '  -2.43%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '0.533'
$ws.Range("E9").Value = '  +1.43%  '
$ws.Range("D10").Value = '0.163'
$ws.Range("E10").Value = '  +2.83%  '
$ws.Range("D11").Value = '6.16'
$ws.Range("E11").Value = '  -2.01%  '
$ws.Range("D12").Value = '0.458'
$ws.Range("E12").Value = '  -1.82%  '
$ws.Range("D13").Value = '37.58'
$ws.Range("E13").Value = '  -1.65%  '
$ws.Range("D14").Value = '0.0000241'
$ws.Range("E14").Value = '  -0.71%  '
$ws.Range("D15").Value = '4.315.98'
$ws.Range("E15").Value = '  -2.52%  '
$ws.Range("D16").Value = '3.700.70'
$ws.Range("E16").Value = '  -2.70%  '
$ws.Range("D17").Value = '67.169.83'
$ws.Range("E17").Value = '  -1.71%  '
$ws.Range("D18").Value = '7.23'
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("E19").Value = '  -1.55%  '
$ws.Range("D20").Value = '17.28'
$ws.Range("E20").Value = '  +8.20%  '
$ws.Range("D21").Value = '484.51'
$ws.Range("E21").Value = '  -0.76%  '
$ws.Range("D22").Value = '9.20'
$ws.Range("E22").Value = '  -1.96%  '
$ws.Range("D23").Value = '0.722'
$ws.Range("E23").Value = '  -0.89%  '
$ws.Range("D24").Value = '84.67'
$ws.Range("E24").Value = '  -1.55%  '
$ws.Range("D25").Value = '0.0000141'
$ws.Range("E25").Value = '  +3.58%  '
$ws.Range("D26").Value = '2.27'
$ws.Range("E26").Value = '  -3.88%  '
$ws.Range("D27").Value = '12.19'
$ws.Range("E27").Value = '  -0.20%  '
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = '9.98'
$ws.Range("E29").Value = '  -1.46%  '
$ws.Range("D30").Value = '2.91'
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("D31").Value = '2.34'
$ws.Range("E31").Value = '  -3.89%  '
$ws.Range("D32").Value = '7.63'
$ws.Range("E32").Value = '  +0.04%  '
$ws.Range("D33").Value = '31.05'
$ws.Range("E33").Value = '  -3.38%  '
$ws.Range("D34").Value = '3.841.00'
$ws.Range("E34").Value = '  -2.58%  '
$ws.Range("E35").Value = '  -2.28%  '
$ws.Range("D36").Value = '3.643.24'
$ws.Range("E36").Value = '  -2.00%  '
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.22%  '
$ws.Range("D38").Value = '0.995'
$ws.Range("E38").Value = '  -1.88%  '
$ws.Range("D39").Value = '5.79'
$ws.Range("E39").Value = '  -0.64%  '
$ws.Range("D40").Value = '0.131'
$ws.Range("E40").Value = '  -2.23%  '
$ws.Range("D41").Value = '0.319'
$ws.Range("E41").Value = '  -1.19%  '
$ws.Range("D42").Value = '48.64'
$ws.Range("E42").Value = '  -1.00%  '
$ws.Range("D43").Value = '423.48'
$ws.Range("E43").Value = '  -5.62%  '
$ws.Range("E44").Value = '  -3.50%  '
$ws.Range("D45").Value = '2.81'
$ws.Range("E45").Value = '  -0.93%  '
$ws.Range("D46").Value = '8.39'
$ws.Range("E46").Value = '  +0.94%  '
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("D48").Value = '40.08'
$ws.Range("E48").Value = '  -3.37%  '
$ws.Range("D49").Value = '140.61'
$ws.Range("E49").Value = '  +1.83%  '
$ws.Range("D50").Value = '2.744.78'
$ws.Range("E50").Value = '  -3.64%  '
$ws.Range("D51").Value = '0.0348'
$ws.Range("E51").Value = '  -0.60%  '
